$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) format, used to restore style
# after forcing a cell to Text format so numeric-looking strings are not
# auto-converted to numbers.
$defaultStyle = $ws.Range("A1").Style

$ws.Range("D2").Value = '66.276.55'
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").Value = '3.567.27'
$ws.Range("E3").Value = '  +5.10%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.64'
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  +2.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.32'
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  +2.86%  '
$ws.Range("D7").Value = '3.565.42'
$ws.Range("E7").Value = '  +5.11%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("E9").Value = '  +3.76%  '
$ws.Range("E10").Value = '  +2.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '8.02'
$ws.Range("D11").Style = $defaultStyle
$ws.Range("E11").Value = '  +1.63%  '
$ws.Range("E12").Value = '  +1.53%  '
$ws.Range("D13").Value = '4.176.21'
$ws.Range("E13").Value = '  +4.95%  '
$ws.Range("E14").Value = '  +4.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '30.18'
$ws.Range("D15").Style = $defaultStyle
$ws.Range("E15").Value = '  +1.49%  '
$ws.Range("D16").Value = '3.567.16'
$ws.Range("E16").Value = '  +4.63%  '
$ws.Range("D17").Value = '66.393.10'
$ws.Range("E17").Value = '  +1.71%  '
$ws.Range("E18").Value = '  -0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.51'
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '  +10.90%  '
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.96'
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '  +1.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '431.73'
$ws.Range("D22").Style = $defaultStyle
$ws.Range("E22").Value = '  +3.90%  '
$ws.Range("E23").Value = '  +5.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.77'
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = '  +1.74%  '
$ws.Range("D25").Value = '3.711.20'
$ws.Range("E25").Value = '  +4.88%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +8.61%  '
$ws.Range("E28").Value = '  +4.30%  '
$ws.Range("E29").Value = '  +2.97%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.19'
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("E32").Value = '  +1.26%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").Value = '3.563.93'
$ws.Range("E34").Value = '  +4.71%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '25.45'
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '  +4.30%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.77'
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  +4.54%  '
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.91'
$ws.Range("D38").Style = $defaultStyle
$ws.Range("E38").Value = '  +5.00%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.65'
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = '  +2.24%  '
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '171.36'
$ws.Range("D41").Style = $defaultStyle
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("E43").Value = '  +3.39%  '
$ws.Range("E44").Value = '  +3.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.96'
$ws.Range("D45").Style = $defaultStyle
$ws.Range("E45").Value = '  +2.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.08'
$ws.Range("D46").Style = $defaultStyle
$ws.Range("E46").Value = '  +1.46%  '
$ws.Range("E47").Value = '  +4.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '25.98'
$ws.Range("D48").Style = $defaultStyle
$ws.Range("E48").Value = '  -2.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.39'
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = '  +4.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.15'
$ws.Range("D50").Style = $defaultStyle
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.38'
$ws.Range("D51").Style = $defaultStyle
$ws.Range("E51").Value = '  +13.70%  '
